$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and week-of dates) ---
$ws.Range("A8").Value = "Volume 31   Number  6"
$ws.Range("C9").Value = "Report Covering the Week  2/5/2024  Through  2/11/2024"

# --- Fix cell formatting (style) where the cell switches between text and numeric ---
$ws.Range("J14").Copy() | Out-Null
$ws.Range("F14","I14","D30","G30","J30").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("C15").Copy() | Out-Null
$ws.Range("G14","H14","C23","C26","D26","E26","C28","C29","C30").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("L15").Copy() | Out-Null
$ws.Range("E30","H30","K30","L30").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# --- Set cell values ---
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = "0"
$ws.Range("H14").Value = "***.*"
$ws.Range("I14").Value = 1
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = -50
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = -33.333333333333
$ws.Range("M15").Value = -42.857142857142
$ws.Range("N15").Value = -75
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -40
$ws.Range("G16").Value = 28
$ws.Range("H16").Value = -28.571428571428
$ws.Range("I16").Value = 27
$ws.Range("J16").Value = 46
$ws.Range("K16").Value = -41.304347826087
$ws.Range("L16").Value = 17.391304347826
$ws.Range("M16").Value = -32.5
$ws.Range("N16").Value = -85.561497326203
$ws.Range("C17").Value = 17
$ws.Range("E17").Value = -5.555555555555
$ws.Range("F17").Value = 60
$ws.Range("G17").Value = 72
$ws.Range("H17").Value = -16.666666666666
$ws.Range("I17").Value = 85
$ws.Range("J17").Value = 104
$ws.Range("K17").Value = -18.269230769230
$ws.Range("L17").Value = 7.594936708860
$ws.Range("M17").Value = 77.083333333333
$ws.Range("N17").Value = -22.018348623853
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -62.5
$ws.Range("I18").Value = 23
$ws.Range("J18").Value = 43
$ws.Range("K18").Value = -46.511627906976
$ws.Range("L18").Value = 27.777777777777
$ws.Range("M18").Value = -59.649122807017
$ws.Range("N18").Value = -93.175074183976
$ws.Range("C19").Value = 29
$ws.Range("D19").Value = 27
$ws.Range("E19").Value = 7.407407407407
$ws.Range("F19").Value = 119
$ws.Range("G19").Value = 110
$ws.Range("H19").Value = 8.181818181818
$ws.Range("I19").Value = 167
$ws.Range("J19").Value = 176
$ws.Range("K19").Value = -5.113636363636
$ws.Range("L19").Value = -5.113636363636
$ws.Range("M19").Value = 63.725490196078
$ws.Range("N19").Value = -5.649717514124
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -20
$ws.Range("G20").Value = 25
$ws.Range("H20").Value = -20
$ws.Range("I20").Value = 28
$ws.Range("J20").Value = 43
$ws.Range("K20").Value = -34.883720930232
$ws.Range("L20").Value = -36.363636363636
$ws.Range("M20").Value = -28.205128205128
$ws.Range("N20").Value = -95.379537953795
$ws.Range("C21").Value = 56
$ws.Range("D21").Value = 63
$ws.Range("E21").Value = -11.111111111111
$ws.Range("F21").Value = 238
$ws.Range("G21").Value = 262
$ws.Range("H21").Value = -9.160305343511
$ws.Range("I21").Value = 335
$ws.Range("J21").Value = 416
$ws.Range("K21").Value = -19.471153846153
$ws.Range("L21").Value = -2.898550724637
$ws.Range("M21").Value = 13.945578231292
$ws.Range("N21").Value = -76.638772663877
$ws.Range("C23").Value = "0"
$ws.Range("D23").Value = 6
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 16
$ws.Range("H23").Value = -68.75
$ws.Range("J23").Value = 23
$ws.Range("K23").Value = -65.217391304347
$ws.Range("L23").Value = -33.333333333333
$ws.Range("M23").Value = 60
$ws.Range("C24").Value = 85
$ws.Range("D24").Value = 79
$ws.Range("E24").Value = 7.594936708860
$ws.Range("F24").Value = 316
$ws.Range("G24").Value = 302
$ws.Range("H24").Value = 4.635761589403
$ws.Range("I24").Value = 490
$ws.Range("J24").Value = 464
$ws.Range("K24").Value = 5.603448275862
$ws.Range("L24").Value = 44.117647058823
$ws.Range("M24").Value = 20.689655172413
$ws.Range("D25").Value = 32
$ws.Range("E25").Value = -12.5
$ws.Range("F25").Value = 111
$ws.Range("G25").Value = 127
$ws.Range("H25").Value = -12.598425196850
$ws.Range("I25").Value = 163
$ws.Range("J25").Value = 185
$ws.Range("K25").Value = -11.891891891891
$ws.Range("L25").Value = 7.947019867549
$ws.Range("M25").Value = -16.836734693877
$ws.Range("C26").Value = "0"
$ws.Range("D26").Value = "0"
$ws.Range("E26").Value = "***.*"
$ws.Range("F26").Value = 9
$ws.Range("H26").Value = 50
$ws.Range("C27").Value = 7
$ws.Range("D27").Value = 4
$ws.Range("E27").Value = 75
$ws.Range("F27").Value = 16
$ws.Range("G27").Value = 13
$ws.Range("H27").Value = 23.076923076923
$ws.Range("I27").Value = 22
$ws.Range("J27").Value = 20
$ws.Range("K27").Value = 10
$ws.Range("L27").Value = 22.222222222222
$ws.Range("C28").Value = "0"
$ws.Range("M28").Value = -25
$ws.Range("C29").Value = "0"
$ws.Range("M29").Value = -25
$ws.Range("C30").Value = "0"
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = -100
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 300
$ws.Range("J30").Value = 1
$ws.Range("K30").Value = 400
$ws.Range("L30").Value = 400
